$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $true, $false, $false, $false, $true, 1, $false, $replace, 2)
}

Replace-Text "M2DocEvaluator.caseUserDoc(M2DocEvaluator.java:1074)" "M2DocEvaluator.caseUserDoc(M2DocEvaluator.java:1120)"
Replace-Text "M2DocEvaluator.doSwitch(M2DocEvaluator.java:1038)" "M2DocEvaluator.doSwitch(M2DocEvaluator.java:1084)"
Replace-Text "M2DocEvaluator.caseBlock(M2DocEvaluator.java:1254)" "M2DocEvaluator.caseBlock(M2DocEvaluator.java:1300)"
Replace-Text "M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:275)" "M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:278)"
Replace-Text "M2DocEvaluator.generate(M2DocEvaluator.java:264)" "M2DocEvaluator.generate(M2DocEvaluator.java:267)"
Replace-Text "M2DocUtils.generate(M2DocUtils.java:712)" "M2DocUtils.generate(M2DocUtils.java:694)"
Replace-Text "AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:459)" "AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:475)"
Replace-Text "AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:369)" "AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:384)"

Write-Output "Done"
